$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header (match formatting of existing header cells: bold + centered)
$ws.Range("R1").Value = "avg_distance_km"
$ws.Range("R1").Font.Bold = $true
$ws.Range("R1").HorizontalAlignment = -4108

# Add the avg_distance_km values for each row
$values = @{
    2  = 4.136376498618044
    3  = 4.171829734708102
    4  = 4.490206231553414
    5  = 4.518039712920642
    6  = 4.508389709238768
    7  = 4.561800693261488
    8  = 4.723351366358533
    9  = 4.541151156676015
    10 = 4.286881911154486
    11 = 4.675176933039502
    12 = 4.299951632670066
    13 = 3.709505535754547
    14 = 3.876734017006266
    15 = 4.476229029663664
    16 = 4.11879765724225
    17 = 4.385888662670677
    18 = 3.76326812791891
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 18).Value = $values[$row]
}
